# Consolidate the workbook down to a single "Prestel" sheet (the old
# Sheet2/Sheet3/Sheet4 scratch-calculation sheets are no longer needed),
# and leave the view focused on the surviving sheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "Prestel"

[void]$wb.Worksheets.Item("Sheet2").Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()
[void]$wb.Worksheets.Item("Sheet4").Delete()

[void]$ws1.Activate()
[void]$ws1.Range("I9").Select()
